$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 35: Graphs / Medium / 417. Pacific Atlantic Water Flow (!) / notes
$ws.Range("A35").Value = "Graphs"
$ws.Range("B35").Value = "Medium"
$ws.Range("C35").Value = "417. Pacific Atlantic Water Flow (!)"

$note = "The core idea is that, instead of going from each cell and checking all sides recursively with dfs to see if it can reach both the oceans,`nWe start our search from the boundaries (top row, bot row, left col, right col), where its guarenteed its touching the ocean, from that we run the dfs. And we track each cell with its corresponding ocean set, so theres no redundant traversal, and we only proceed with the dfs if its a valid cell ie. index within bounds, prev height < curr height, curr not in set"
$ws.Range("D35").Value = $note

# Match formatting used by the other "Graphs" rows (B/C = Neutral style, D = wrap/top alignment)
$ws.Range("B35").Style = "Neutral"
$ws.Range("C35").Style = "Neutral"
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").VerticalAlignment = -4160
$ws.Range("D35").WrapText = $true

$ws.Rows.Item(35).RowHeight = 60.6

# Hyperlink for the new problem name cell
$url = "https://leetcode.com/problems/pacific-atlantic-water-flow/"
$ws.Hyperlinks.Add($ws.Range("C35"), $url, [Type]::Missing, [Type]::Missing, $url)

# Update the active selection to match the saved workbook state
$ws.Range("B38").Select()
